# Auto-generated: update crypto price (D) and 1h volume change (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.604.57'
$ws.Range("E2").Value = '  -1.00%  '
$ws.Range("D3").Value = '2.644.85'
$ws.Range("E3").Value = '  +1.26%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '518.13'
$ws.Range("E5").Value = '  -0.63%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.74'
$ws.Range("E6").Value = '  -1.42%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.996'
$ws.Range("E7").Value = '  -0.32%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.573'
$ws.Range("E8").Value = '  +0.61%  '
$ws.Range("D9").Value = '2.670.05'
$ws.Range("E9").Value = '  +2.00%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.50'
$ws.Range("E10").Value = '  +2.93%  '
$ws.Range("E11").Value = '  +1.70%  '
$ws.Range("E12").Value = '  -0.24%  '
$ws.Range("E13").Value = '  -1.31%  '
$ws.Range("D14").Value = '3.106.33'
$ws.Range("E14").Value = '  +1.26%  '
$ws.Range("D15").Value = '59.471.79'
$ws.Range("E15").Value = '  -1.24%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.28'
$ws.Range("E16").Value = '  -0.22%  '
$ws.Range("E17").Value = '  +0.68%  '
$ws.Range("D18").Value = '2.656.18'
$ws.Range("E18").Value = '  +1.61%  '
$ws.Range("E19").Value = '  -0.23%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '346.87'
$ws.Range("E21").Value = '  +1.06%  '
$ws.Range("E22").Value = '  +1.34%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("E23").Value = '  +0.42%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.93'
$ws.Range("E24").Value = '  +1.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.424'
$ws.Range("E25").Value = '  +1.37%  '
$ws.Range("D26").Value = '2.764.62'
$ws.Range("E26").Value = '  +1.31%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.998'
$ws.Range("E27").Value = '  -0.03%  '
$ws.Range("E28").Value = '  -0.93%  '
$ws.Range("E29").Value = '  +2.44%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.24'
$ws.Range("E30").Value = '  +2.40%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.998'
$ws.Range("E31").Value = '  -0.23%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.53'
$ws.Range("E32").Value = '  +8.25%  '
$ws.Range("E33").Value = '  +0.44%  '
$ws.Range("E34").Value = '  -0.49%  '
$ws.Range("E35").Value = '  +16.84%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '149.63'
$ws.Range("E36").Value = '  -0.11%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.04'
$ws.Range("E37").Value = '  +2.11%  '
$ws.Range("E38").Value = '  +2.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.874'
$ws.Range("E39").Value = '  +1.11%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.76'
$ws.Range("E40").Value = '  +1.10%  '
$ws.Range("E41").Value = '  +3.41%  '
$ws.Range("E42").Value = '  -0.58%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '284.60'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.618'
$ws.Range("E44").Value = '  -0.99%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0996'
$ws.Range("E45").Value = '  -1.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.993'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '19.83'
$ws.Range("E47").Value = '  +1.51%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0546'
$ws.Range("E48").Value = '  -0.24%  '
$ws.Range("E49").Value = '  +0.26%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.76'
$ws.Range("E50").Value = '  +1.84%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '10.27'
$ws.Range("E51").Value = '  -1.27%  '
